$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (ConvexHullArea) holds computed values with many decimal places.
# Round each value in D2:D117 to the nearest whole number, replacing the
# formula-free numeric literal in place (matching the commit's CSV export
# behavior of truncating precision to integers).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($null -ne $val) {
        $cell.Value = [Math]::Floor([double]$val + 0.5)
    }
}
